# 11. Hafta Bütçe Güncellemesi
# 11. Hafta için bütçe güncellemesi yapılmıştır.

$wb = $excel.ActiveWorkbook

# --- Delete the obsolete "BBM487-20142-1 Bütçe" sheet (its data now lives on the
#     "Bilgilendirme" sheet's own budget table) ---
$budgetSheet = $wb.Worksheets.Item("BBM487-20142-1 Bütçe")
$budgetSheet.Delete() | Out-Null

# --- Work on the main informational sheet ---
$ws = $wb.Worksheets.Item("BBM487-20142-1 Bilgilendirme")
$ws.Activate() | Out-Null

# Bump the "Tarih" reference date by one week (1 May -> 8 May 2015)
$ws.Range("N9").Formula = "=DATE(2015,5,8)"

# Fill in week 11 (row 42) of the budget table with this week's real figures
$ws.Range("C42").Value = 10000
$ws.Range("D42").Value = 6000
$ws.Range("E42").Value = 6000
$ws.Range("F42").Value = 4000
$ws.Range("G42").Value = 4000
$ws.Range("H42").Value = 4000
$ws.Range("I42").Value = 0
$ws.Range("J42").Formula = "=SUM(C42:I42)"
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("N42").Formula = "= ( (M42 * L42) * K42 / 100 ) + K42"
$ws.Range("O42").Formula = "=N42 - J42"
$ws.Range("P42").Formula = "=Q42 * 0.1"
$ws.Range("Q42").Formula = "= (Q41 + O42) + P41"

# Totals row now needs to span through the newly-filled row 42
$ws.Range("J45").Formula = "=SUM(J32:J42)"
$ws.Range("N45").Formula = "=SUM(N32:N42)"
$ws.Range("P45").Formula = "=SUM(P32:P44) - P342"
$ws.Range("Q45").Formula = "= Q42"

# Scroll/selection bookkeeping to match where the editor left off
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("N10").Select() | Out-Null
